# Update tab names in template4.xlsx
#  - ETPT_ATT_JUR      -> ETPT_ATTJ      (also un-hide it, it becomes the active tab)
#  - ETPT_ATT_JUR_DDG  -> ETPT_ATTJ_DDG
#  - Update selections on ETPT_ATTJ_DDG and ETPT_TJ_DDG to I5
# Renaming the sheets automatically re-points every formula that referenced
# the old sheet names (e.g. ETPT_ATT_JUR_DDG!D5 -> ETPT_ATTJ_DDG!D5).

$wb = $excel.ActiveWorkbook

$wsAttJur = $wb.Worksheets.Item("ETPT_ATT_JUR")
$wsAttJur.Name = "ETPT_ATTJ"

$wsAttJurDdg = $wb.Worksheets.Item("ETPT_ATT_JUR_DDG")
$wsAttJurDdg.Name = "ETPT_ATTJ_DDG"

# ETPT_ATTJ is no longer hidden and becomes the selected/active tab.
$wsAttJur.Visible = -1

# Selection moves to I5 on these two sheets.
$null = $wsAttJurDdg.Range("I5").Select()

$wsTjDdg = $wb.Worksheets.Item("ETPT_TJ_DDG")
$null = $wsTjDdg.Range("I5").Select()

# Make ETPT_ATTJ the active sheet/tab (was ETPT_ATT_JUR_DDG before).
$wsAttJur.Activate()
